$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.434.42'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '3.427.20'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.26'
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '621.43'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.42'
$ws.Range('E7').Value = '  +1.38%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.396'
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.976'
$ws.Range('E10').Value = '  +3.29%  '
$ws.Range('D11').Value = '3.425.23'
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.15'
$ws.Range('E12').Value = '  +6.30%  '
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.32'
$ws.Range('E14').Value = '  +4.53%  '
$ws.Range('D15').Value = '93.303.54'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '4.064.02'
$ws.Range('E16').Value = '  +1.65%  '
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.23'
$ws.Range('E18').Value = '  +2.64%  '
$ws.Range('D19').Value = '3.431.75'
$ws.Range('E19').Value = '  +1.85%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.18'
$ws.Range('E20').Value = '  +7.03%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.70'
$ws.Range('E21').Value = '  +4.98%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.39'
$ws.Range('E22').Value = '  +7.69%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '502.10'
$ws.Range('E23').Value = '  +2.99%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.471'
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('E25').Value = '  +8.10%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000185'
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '95.21'
$ws.Range('E27').Value = '  +5.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.07'
$ws.Range('E28').Value = '  +4.53%  '
$ws.Range('D29').Value = '3.611.17'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '11.40'
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.139'
$ws.Range('E32').Value = '  +5.08%  '
$ws.Range('E33').Value = '  +3.70%  '
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.173'
$ws.Range('E35').Value = '  +1.17%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '29.62'
$ws.Range('E36').Value = '  +3.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.556'
$ws.Range('E37').Value = '  +4.61%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '562.39'
$ws.Range('E38').Value = '  +5.59%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.51'
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.41'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.909'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.72'
$ws.Range('E44').Value = '  +3.45%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '23.66'
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.70'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0413'
$ws.Range('E47').Value = '  +4.53%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.51'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '53.46'
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.14'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.15'
$ws.Range('E51').Value = '  +3.86%  '
